# Auto-generated script to update Moogle_Profits market price data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 876
$ws.Range("I2").Value = 307.75
$ws.Range("J2").Value = 1850.1428
$ws.Range("K2").Value = 307.75
$ws.Range("L2").Value = 1850.1428
$ws.Range("M2").Value = -194.75
$ws.Range("N2").Value = -2076.1428

$ws.Range("H17").Value = 2461.2
$ws.Range("J17").Value = 2461.2
$ws.Range("L17").Value = 7383.599999999999
$ws.Range("N17").Value = -7719.599999999999

$ws.Range("H32").Value = 4634.909
$ws.Range("I32").Value = 3996.25
$ws.Range("K32").Value = 3996.25
$ws.Range("M32").Value = -3670.25

$ws.Range("H113").Value = 3671.125
$ws.Range("I113").Value = 3053.8
$ws.Range("J113").Value = 4700
$ws.Range("K113").Value = 3053.8
$ws.Range("L113").Value = 4700
$ws.Range("M113").Value = 200.1999999999998
$ws.Range("N113").Value = -11208

$ws.Range("H135").Value = 807.62067
$ws.Range("I135").Value = 612.38464
$ws.Range("J135").Value = 2499.6667
$ws.Range("K135").Value = 5511.46176
$ws.Range("L135").Value = 22497.0003
$ws.Range("M135").Value = -2976.46176
$ws.Range("N135").Value = -27567.0003

$ws.Range("H138").Value = 4945.227
$ws.Range("I138").Value = 4075.725
$ws.Range("K138").Value = 12227.175
$ws.Range("M138").Value = -7087.174999999999

$ws.Range("H141").Value = 3140.1428
$ws.Range("I141").Value = 1102.8422
$ws.Range("K141").Value = 3308.5266
$ws.Range("M141").Value = 1871.4734


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6996.3584
$ws.Range("I32").Value = 3030.0168
$ws.Range("J32").Value = 36248.125
$ws.Range("K32").Value = 3030.0168
$ws.Range("L32").Value = 36248.125
$ws.Range("M32").Value = -2743.0168
$ws.Range("N32").Value = -36822.125

$ws.Range("H61").Value = 4801.1577
$ws.Range("I61").Value = 4990.1113
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 4990.1113
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -4778.1113
$ws.Range("N61").Value = -1824

$ws.Range("H102").Value = 968
$ws.Range("I102").Value = 968
$ws.Range("K102").Value = 968
$ws.Range("M102").Value = 654

$ws.Range("H131").Value = 163000
$ws.Range("J131").Value = 163000
$ws.Range("L131").Value = 163000
$ws.Range("N131").Value = -173080

$ws.Range("H132").Value = 3677.8
$ws.Range("I132").Value = 1947.9062
$ws.Range("J132").Value = 10597.375
$ws.Range("K132").Value = 5843.7186
$ws.Range("L132").Value = 31792.125
$ws.Range("M132").Value = -3313.7186
$ws.Range("N132").Value = -36852.125

$ws.Range("H136").Value = 4801.1577
$ws.Range("I136").Value = 4990.1113
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 14970.3339
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -12420.3339
$ws.Range("N136").Value = -9300


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3896.4375
$ws.Range("I86").Value = 4395.5454
$ws.Range("J86").Value = 2798.4
$ws.Range("K86").Value = 4395.5454
$ws.Range("L86").Value = 2798.4
$ws.Range("M86").Value = -3272.5454
$ws.Range("N86").Value = -5044.4

$ws.Range("H89").Value = 3896.4375
$ws.Range("I89").Value = 4395.5454
$ws.Range("J89").Value = 2798.4
$ws.Range("K89").Value = 21977.727
$ws.Range("L89").Value = 13992
$ws.Range("M89").Value = -16361.727
$ws.Range("N89").Value = -25224

$ws.Range("H105").Value = 4025.5293
$ws.Range("I105").Value = 3922.3572
$ws.Range("K105").Value = 3922.3572
$ws.Range("M105").Value = -2175.3572

$ws.Range("H134").Value = 2350.7
$ws.Range("I134").Value = 1696.6786
$ws.Range("K134").Value = 5090.0358
$ws.Range("M134").Value = -2555.0358


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2460.6
$ws.Range("J22").Value = 3794.6667
$ws.Range("L22").Value = 3794.6667
$ws.Range("N22").Value = -4494.6667

$ws.Range("H31").Value = 8584.23
$ws.Range("I31").Value = 2914.6
$ws.Range("K31").Value = 2914.6
$ws.Range("M31").Value = -2619.6

$ws.Range("H34").Value = 8584.23
$ws.Range("I34").Value = 2914.6
$ws.Range("K34").Value = 2914.6
$ws.Range("M34").Value = -2712.6

$ws.Range("H99").Value = 2008.2413
$ws.Range("I99").Value = 2032.3846
$ws.Range("K99").Value = 2032.3846
$ws.Range("M99").Value = -534.3846000000001

$ws.Range("H126").Value = 2008.2413
$ws.Range("I126").Value = 2032.3846
$ws.Range("K126").Value = 6097.1538
$ws.Range("M126").Value = -3627.1538


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 495.7143
$ws.Range("I35").Value = 343.33334
$ws.Range("J35").Value = 610
$ws.Range("K35").Value = 1030.00002
$ws.Range("L35").Value = 1830
$ws.Range("M35").Value = -742.0000199999999
$ws.Range("N35").Value = -2406

$ws.Range("H68").Value = 796.6667
$ws.Range("I68").Value = 594.625
$ws.Range("K68").Value = 1783.875
$ws.Range("M68").Value = -972.875

$ws.Range("H71").Value = 796.6667
$ws.Range("I71").Value = 594.625
$ws.Range("K71").Value = 5351.625
$ws.Range("M71").Value = -1295.625

$ws.Range("H88").Value = 17999.8
$ws.Range("I88").Value = 9999
$ws.Range("J88").Value = 20000
$ws.Range("K88").Value = 29997
$ws.Range("L88").Value = 60000
$ws.Range("M88").Value = -29569
$ws.Range("N88").Value = -60856

$ws.Range("H91").Value = 17999.8
$ws.Range("I91").Value = 9999
$ws.Range("J91").Value = 20000
$ws.Range("K91").Value = 29997
$ws.Range("L91").Value = 60000
$ws.Range("M91").Value = -28515
$ws.Range("N91").Value = -62964

$ws.Range("H94").Value = 11199
$ws.Range("I94").Value = 3497.5
$ws.Range("J94").Value = 16333.333
$ws.Range("K94").Value = 10492.5
$ws.Range("L94").Value = 48999.999
$ws.Range("M94").Value = -9816.5
$ws.Range("N94").Value = -50351.999

$ws.Range("H97").Value = 1204.5555
$ws.Range("I97").Value = 1177.3334
$ws.Range("J97").Value = 1259
$ws.Range("K97").Value = 3532.0002
$ws.Range("L97").Value = 3777
$ws.Range("M97").Value = -3036.0002
$ws.Range("N97").Value = -4769

$ws.Range("H107").Value = 434.36365
$ws.Range("J107").Value = 471.22223
$ws.Range("L107").Value = 1413.66669
$ws.Range("N107").Value = -5253.66669

$ws.Range("H140").Value = 1938.7906
$ws.Range("I140").Value = 1758.4
$ws.Range("K140").Value = 5275.200000000001
$ws.Range("M140").Value = -95.20000000000073


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7233.375
$ws.Range("I70").Value = 5749
$ws.Range("K70").Value = 5749
$ws.Range("M70").Value = -5479

$ws.Range("H73").Value = 7233.375
$ws.Range("I73").Value = 5749
$ws.Range("K73").Value = 5749
$ws.Range("M73").Value = -4813

$ws.Range("H132").Value = 2887.7036
$ws.Range("I132").Value = 1998.8334
$ws.Range("J132").Value = 9998.666999999999
$ws.Range("K132").Value = 5996.5002
$ws.Range("L132").Value = 29996.001
$ws.Range("M132").Value = -3466.5002
$ws.Range("N132").Value = -35056.001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1639.2941
$ws.Range("I22").Value = 1584.7142
$ws.Range("J22").Value = 1677.5
$ws.Range("K22").Value = 1584.7142
$ws.Range("L22").Value = 1677.5
$ws.Range("M22").Value = -1289.7142
$ws.Range("N22").Value = -2267.5

$ws.Range("H27").Value = 1639.2941
$ws.Range("I27").Value = 1584.7142
$ws.Range("J27").Value = 1677.5
$ws.Range("K27").Value = 1584.7142
$ws.Range("L27").Value = 1677.5
$ws.Range("M27").Value = -1477.7142
$ws.Range("N27").Value = -1891.5

$ws.Range("H40").Value = 6508.905
$ws.Range("I40").Value = 5523.706
$ws.Range("J40").Value = 10696
$ws.Range("K40").Value = 5523.706
$ws.Range("L40").Value = 10696
$ws.Range("M40").Value = -5387.706
$ws.Range("N40").Value = -10968

$ws.Range("H122").Value = 4167.7393
$ws.Range("I122").Value = 4701.1113
$ws.Range("K122").Value = 14103.3339
$ws.Range("M122").Value = -11653.3339

$ws.Range("H132").Value = 2642.7026
$ws.Range("I132").Value = 1659.32
$ws.Range("J132").Value = 4691.4165
$ws.Range("K132").Value = 4977.96
$ws.Range("L132").Value = 14074.2495
$ws.Range("M132").Value = -2447.96
$ws.Range("N132").Value = -19134.2495

$ws.Range("H136").Value = 5610.297
$ws.Range("I136").Value = 2934.25
$ws.Range("J136").Value = 7691.6665
$ws.Range("K136").Value = 8802.75
$ws.Range("L136").Value = 23074.9995
$ws.Range("M136").Value = -6252.75
$ws.Range("N136").Value = -28174.9995


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 38995
$ws.Range("J69").Value = 38995
$ws.Range("L69").Value = 38995
$ws.Range("N69").Value = -40493

$ws.Range("H72").Value = 38995
$ws.Range("J72").Value = 38995
$ws.Range("L72").Value = 116985
$ws.Range("N72").Value = -124473

$ws.Range("H81").Value = 1549.1428
$ws.Range("J81").Value = 1265.8334
$ws.Range("L81").Value = 2531.6668
$ws.Range("N81").Value = -4653.6668

$ws.Range("H84").Value = 1549.1428
$ws.Range("J84").Value = 1265.8334
$ws.Range("L84").Value = 12658.334
$ws.Range("N84").Value = -23266.334

$ws.Range("H96").Value = 5052.4736
$ws.Range("I96").Value = 2692.4285
$ws.Range("J96").Value = 6429.1665
$ws.Range("K96").Value = 2692.4285
$ws.Range("L96").Value = 6429.1665
$ws.Range("M96").Value = -1319.4285
$ws.Range("N96").Value = -9175.166499999999

$ws.Range("H122").Value = 1963.4546
$ws.Range("I122").Value = 1973.5897
$ws.Range("J122").Value = 1884.4
$ws.Range("K122").Value = 5920.7691
$ws.Range("L122").Value = 5653.200000000001
$ws.Range("M122").Value = -3470.7691
$ws.Range("N122").Value = -10553.2

$ws.Range("H132").Value = 3771.9285
$ws.Range("I132").Value = 2510.7727
$ws.Range("K132").Value = 7532.3181
$ws.Range("M132").Value = -5002.3181

$ws.Range("H136").Value = 2597.4375
$ws.Range("J136").Value = 6043.1665
$ws.Range("L136").Value = 18129.4995
$ws.Range("N136").Value = -23229.4995

